$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed market data (price + 1h volume change) for each coin row.
# Cells whose new text would otherwise be auto-coerced to a number by Excel
# are forced back to text ('@') first, matching the source column's text storage.

$ws.Range('D2').Value = '26.716.17'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.638.36'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  -0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.76'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  -0.64%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0622'
$ws.Range('E9').Value = '  -0.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.07'
$ws.Range('E10').Value = '  -0.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('D12').Value = '1.865.98'
$ws.Range('E12').Value = '  -0.71%  '
$ws.Range('D13').Value = '1.623.46'
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('E14').Value = '  -1.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  -1.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.45'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').Value = '26.688.21'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('E18').Value = '  -2.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '211.20'
$ws.Range('E19').Value = '  -3.33%  '
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('E21').Value = '  -0.68%  '
$ws.Range('E22').Value = '  -1.25%  '
$ws.Range('E23').Value = '  -3.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.25'
$ws.Range('E24').Value = '  -2.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.66'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('E26').Value = '  -0.41%  '
$ws.Range('E27').Value = '  -1.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.07'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.55'
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('E30').Value = '  -2.93%  '
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('D34').Value = '1.268.89'
$ws.Range('E34').Value = '  -1.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.52'
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('E36').Value = '  -0.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0174'
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.528'
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('B39').Value = 'PaxDollar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.01'
$ws.Range('E39').Value = '  -0.30%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.802'
$ws.Range('E40').Value = '  -3.17%  '
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.19'
$ws.Range('E42').Value = '  -3.36%  '
$ws.Range('E43').Value = '  -3.68%  '
$ws.Range('D44').Value = '1.775.58'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '91.26'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.96'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('E47').Value = '  -1.76%  '
$ws.Range('E48').Value = '  +0.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.55'
$ws.Range('E49').Value = '  -2.66%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0959'
$ws.Range('E50').Value = '  -1.18%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.407'
$ws.Range('E51').Value = '  -0.58%  '
